# "Fruta / hortaliza, semanal" — weekly refresh of the price table.
# A new weekly record is inserted at row 48 (pushing the existing
# rows 48-95 down to 49-96, dimension grows from A1:R95 to A1:R96).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 48, shifting everything
# from row 48 downward (previous row 48 becomes row 49, ..., previous
# row 95 becomes row 96).
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with this week's record.
$ws.Cells.Item(48, 1).Value  = 1
$ws.Cells.Item(48, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value  = 44792
$ws.Cells.Item(48, 5).Value  = 15
$ws.Cells.Item(48, 6).Value  = 100112021
$ws.Cells.Item(48, 7).Value  = "Ají"
$ws.Cells.Item(48, 8).Value  = "Inferno"
$ws.Cells.Item(48, 9).Value  = "Primera"
$ws.Cells.Item(48, 10).Value = 130
$ws.Cells.Item(48, 11).Value = 10000
$ws.Cells.Item(48, 12).Value = 11000
$ws.Cells.Item(48, 13).Value = 10500
$ws.Cells.Item(48, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 700
$ws.Cells.Item(48, 17).Value = 15
$ws.Cells.Item(48, 18).Value = "Hortaliza"
